$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.193.12"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "1.858.03"
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.98"
$ws.Range("E5").Value = "  +3.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.14"
$ws.Range("E8").Value = "  +6.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.329"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0692"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("D12").Value = "2.126.88"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.48"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "1.846.21"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.676"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").Value = "35.209.71"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.85"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "0.0₃0795"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.39"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.20"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.75"
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.30"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("E26").Value = "  +26.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.00"
$ws.Range("E27").Value = "  +3.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.64"
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  +1.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.99"
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.83"
$ws.Range("E33").Value = "  +28.53%  "
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("E35").Value = "  +9.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.812"
$ws.Range("E36").Value = "  +16.30%  "
$ws.Range("E37").Value = "  +7.78%  "
$ws.Range("E38").Value = "  +3.98%  "
$ws.Range("E39").Value = "  +4.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "89.91"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D41").Value = "1.345.41"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0595"
$ws.Range("E42").Value = "  +14.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.97"
$ws.Range("E43").Value = "  +3.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.32"
$ws.Range("E44").Value = "  +3.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.42"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.54"
$ws.Range("E46").Value = "  +43.88%  "
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("E48").Value = "  +5.60%  "
$ws.Range("D49").Value = "2.047.29"
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0680"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("E51").Value = "  +0.41%  "
